$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-31 from 45221 to 45224
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
